$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.753.56"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.464.71"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").Value = "3.458.44"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.566"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "4.019.69"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "585.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "3.468.30"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "69.776.71"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.852"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "96.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.78%  "
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("E32").Value = "  -5.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "587.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -15.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0479"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0961"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.62%  "
$ws.Range("D43").Value = "3.257.25"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "30.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
